$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($r=1; $r -le 7; $r++) {
  for ($c=1; $c -le 11; $c++) {
    $cell = $ws.Cells.Item($r,$c)
    $v = $cell.Value2
    Write-Host "R${r}C${c} val=[$v]"
  }
}
